# This script updates the LocationX values (column B) of Sheet1 to reflect
# a recomputed 3D trajectory transformation where predicted Z coordinates
# below 0 have been filtered out before re-projecting to the X location.
#
# The new values below are the recomputed LocationX values for rows 2-62
# (time = 0 .. 2, step 1/30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    838.7376338955329,
    1584.131602906735,
    2495.20737646486,
    2959.437207276416,
    3607.538196324411,
    3976.520384812403,
    4416.985764789569,
    4970.312816123001,
    5348.458202839564,
    5638.494694604529,
    5933.966239424018,
    6346.166754829405,
    6493.441058445035,
    6865.611960649431,
    7103.248667245887,
    7226.344060710124,
    7417.451344920351,
    7708.274839746323,
    7772.717299292321,
    8015.928693941787,
    8254.578600193758,
    8413.208558401786,
    8448.386311720573,
    8661.708609395282,
    8911.114864939369,
    9159.746586652001,
    9401.129076062811,
    9494.715932235424,
    9597.654090265871,
    9735.915919783129,
    10107.16445152251,
    10369.83513935365,
    10497.53168777822,
    10492.64187557545,
    10511.74467251392,
    10659.45536933116,
    10780.6741020643,
    10899.81047866953,
    11057.98904495969,
    11121.69890184989,
    11152.04992793766,
    11313.62404544317,
    11488.55564298573,
    11611.62376677069,
    11658.72075392685,
    11623.79476981653,
    11825.4625923721,
    11613.85853947689,
    11689.23638223392,
    12059.9044354463,
    11870.69058719389,
    11903.98182651385,
    12121.67493443103,
    12196.36721382907,
    12108.72989448313,
    12244.93390605437,
    12229.97123178263,
    12348.4553037889,
    12348.4553037889,
    12348.4553037889,
    12348.4553037889
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}
